$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Septiembre de 2020 a las 07:24"

# Row 5 - India
$ws.Range("B5").Value = 4846427
$ws.Range("C5").Value = 1424
$ws.Range("D5").Value = 3780107
$ws.Range("E5").Value = 986566

# Row 27 - Israel
$ws.Range("B27").Value = 156596
$ws.Range("C27").Value = 992
$ws.Range("D27").Value = 115122
$ws.Range("E27").Value = 40355

# Row 43 - China
$ws.Range("B43").Value = 85194
$ws.Range("C43").Value = 10
$ws.Range("D43").Value = 80415
$ws.Range("E43").Value = 145

# Row 60 - Uzbekistan
$ws.Range("B60").Value = 47423
$ws.Range("C60").Value = 136
$ws.Range("D60").Value = 44002
$ws.Range("E60").Value = 3029
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 392

# Row 64 - Kirguistan
$ws.Range("B64").Value = 44928
$ws.Range("C64").Value = 47
$ws.Range("D64").Value = 41023
$ws.Range("E64").Value = 2842

# Row 174 - Papua Nueva Guinea
$ws.Range("B174").Value = 511
$ws.Range("C174").Value = 1
$ws.Range("E174").Value = 273

# Row 187 - Butan
$ws.Range("B187").Value = 245
$ws.Range("C187").Value = 1
$ws.Range("D187").Value = 161
$ws.Range("E187").Value = 84
